# Scheduled runner update: refresh cached market-board price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit* columns) for the rows whose
# underlying item prices changed, across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 590002.2
$ws.Range("I19").Value = 1667357.6
$ws.Range("K19").Value = 1667357.6
$ws.Range("M19").Value = -1667182.6

$ws.Range("H32").Value = 2777.8
$ws.Range("J32").Value = 2777.8
$ws.Range("L32").Value = 2777.8
$ws.Range("N32").Value = -3429.8

$ws.Range("H42").Value = 525.5
$ws.Range("J42").Value = 904.5
$ws.Range("L42").Value = 2713.5
$ws.Range("N42").Value = -3173.5

$ws.Range("H43").Value = 1416.2307
$ws.Range("I43").Value = 1299.3334
$ws.Range("J43").Value = 1451.3
$ws.Range("K43").Value = 1299.3334
$ws.Range("L43").Value = 1451.3
$ws.Range("M43").Value = -1230.3334
$ws.Range("N43").Value = -1589.3

$ws.Range("H112").Value = 1749.1892
$ws.Range("J112").Value = 1821.8485
$ws.Range("L112").Value = 5465.5455
$ws.Range("N112").Value = -7681.5455

$ws.Range("H135").Value = 705.5789
$ws.Range("I135").Value = 385.16666
$ws.Range("J135").Value = 1254.8572
$ws.Range("K135").Value = 3466.49994
$ws.Range("L135").Value = 11293.7148
$ws.Range("M135").Value = -931.4999399999997
$ws.Range("N135").Value = -16363.7148

$ws.Range("H137").Value = 1651.3
$ws.Range("I137").Value = 1184.6666
$ws.Range("K137").Value = 3553.9998
$ws.Range("M137").Value = -1003.9998

$ws.Range("H141").Value = 2334828.5
$ws.Range("I141").Value = 2546176.5
$ws.Range("J141").Value = 9999
$ws.Range("K141").Value = 7638529.5
$ws.Range("L141").Value = 29997
$ws.Range("M141").Value = -7633349.5
$ws.Range("N141").Value = -40357

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3898.5881
$ws.Range("I32").Value = 3440.1538
$ws.Range("K32").Value = 3440.1538
$ws.Range("M32").Value = -3153.1538

$ws.Range("H61").Value = 6774.7
$ws.Range("J61").Value = 5779.8
$ws.Range("L61").Value = 5779.8
$ws.Range("N61").Value = -6203.8

$ws.Range("H74").Value = 939.6667
$ws.Range("I74").Value = 536.3103599999999
$ws.Range("J74").Value = 3864
$ws.Range("K74").Value = 536.3103599999999
$ws.Range("L74").Value = 3864
$ws.Range("M74").Value = 337.6896400000001
$ws.Range("N74").Value = -5612

$ws.Range("H77").Value = 939.6667
$ws.Range("I77").Value = 536.3103599999999
$ws.Range("J77").Value = 3864
$ws.Range("K77").Value = 2681.5518
$ws.Range("L77").Value = 19320
$ws.Range("M77").Value = 1686.4482
$ws.Range("N77").Value = -28056

$ws.Range("H102").Value = 1208.8334
$ws.Range("I102").Value = 1050.75
$ws.Range("K102").Value = 1050.75
$ws.Range("M102").Value = 571.25

$ws.Range("H136").Value = 6774.7
$ws.Range("J136").Value = 5779.8
$ws.Range("L136").Value = 17339.4
$ws.Range("N136").Value = -22439.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2320.4
$ws.Range("I105").Value = 2170.8518
$ws.Range("K105").Value = 2170.8518
$ws.Range("M105").Value = -423.8517999999999

$ws.Range("H134").Value = 8727
$ws.Range("I134").Value = 10080.667
$ws.Range("J134").Value = 3312.3333
$ws.Range("K134").Value = 30242.001
$ws.Range("L134").Value = 9936.999899999999
$ws.Range("M134").Value = -27707.001
$ws.Range("N134").Value = -15006.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2776.8
$ws.Range("I31").Value = 2475
$ws.Range("J31").Value = 2978
$ws.Range("K31").Value = 2475
$ws.Range("L31").Value = 2978
$ws.Range("M31").Value = -2180
$ws.Range("N31").Value = -3568

$ws.Range("H34").Value = 2776.8
$ws.Range("I34").Value = 2475
$ws.Range("J34").Value = 2978
$ws.Range("K34").Value = 2475
$ws.Range("L34").Value = 2978
$ws.Range("M34").Value = -2273
$ws.Range("N34").Value = -3382

$ws.Range("H58").Value = 1892462.2
$ws.Range("I58").Value = 2072315.9
$ws.Range("K58").Value = 2072315.9
$ws.Range("M58").Value = -2072112.9

$ws.Range("H88").Value = 22500
$ws.Range("J88").Value = 22500
$ws.Range("L88").Value = 22500
$ws.Range("N88").Value = -23312

$ws.Range("H91").Value = 22500
$ws.Range("J91").Value = 22500
$ws.Range("L91").Value = 22500
$ws.Range("N91").Value = -25308

$ws.Range("H95").Value = 26123.6
$ws.Range("J95").Value = 26123.6
$ws.Range("L95").Value = 26123.6
$ws.Range("N95").Value = -31615.6

$ws.Range("H136").Value = 1892462.2
$ws.Range("I136").Value = 2072315.9
$ws.Range("K136").Value = 6216947.699999999
$ws.Range("M136").Value = -6214397.699999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 312
$ws.Range("I6").Value = 312
$ws.Range("K6").Value = 936
$ws.Range("M6").Value = -823

$ws.Range("H11").Value = 718.1667
$ws.Range("I11").Value = 662
$ws.Range("K11").Value = 1986
$ws.Range("M11").Value = -1846

$ws.Range("H64").Value = 2962.5
$ws.Range("J64").Value = 3416.6667
$ws.Range("L64").Value = 10250.0001
$ws.Range("N64").Value = -10790.0001

$ws.Range("H67").Value = 2962.5
$ws.Range("J67").Value = 3416.6667
$ws.Range("L67").Value = 10250.0001
$ws.Range("N67").Value = -12122.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2767.0833
$ws.Range("I80").Value = 2519.4
$ws.Range("J80").Value = 4005.5
$ws.Range("K80").Value = 2519.4
$ws.Range("L80").Value = 4005.5
$ws.Range("M80").Value = -1521.4
$ws.Range("N80").Value = -6001.5

$ws.Range("H83").Value = 2767.0833
$ws.Range("I83").Value = 2519.4
$ws.Range("J83").Value = 4005.5
$ws.Range("K83").Value = 12597
$ws.Range("L83").Value = 20027.5
$ws.Range("M83").Value = -7605
$ws.Range("N83").Value = -30011.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 15000
$ws.Range("J97").Value = 15000
$ws.Range("L97").Value = 15000
$ws.Range("N97").Value = -16982

$ws.Range("H132").Value = 2578.7646
$ws.Range("I132").Value = 1710.3
$ws.Range("K132").Value = 5130.9
$ws.Range("M132").Value = -2600.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1749.5
$ws.Range("I100").Value = 1500
$ws.Range("J100").Value = 1999
$ws.Range("K100").Value = 3000
$ws.Range("L100").Value = 3998
$ws.Range("M100").Value = -2459
